$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 100, shifting rows 100:110 down to 102:112.
$ws.Range("A100:T101").EntireRow.Insert()

# New row 100 data
$ws.Range("A100").Value = 3
$ws.Range("B100").Value = "Femacal de La Calera"
$ws.Range("C100").Value = "Coquimbo"
$ws.Range("D100").Value = 45194
$ws.Range("E100").Value = 5
$ws.Range("F100").Value = "Fruta"
$ws.Range("G100").Value = 100108
$ws.Range("H100").Value = "Tropicales y subtropicales"
$ws.Range("I100").Value = 100108004
$ws.Range("J100").Value = "Papaya"
$ws.Range("K100").Value = "Cultivar IV Región"
$ws.Range("L100").Value = "Primera"
$ws.Range("M100").Value = 56
$ws.Range("N100").Value = 18000
$ws.Range("O100").Value = 18000
$ws.Range("P100").Value = 18000
$ws.Range("Q100").Value = "`$/bandeja 10 kilos"
$ws.Range("R100").Value = "Provincia del Elquí"
$ws.Range("S100").Value = 1800
$ws.Range("T100").Value = 10

# New row 101 data
$ws.Range("A101").Value = 3
$ws.Range("B101").Value = "Femacal de La Calera"
$ws.Range("C101").Value = "Coquimbo"
$ws.Range("D101").Value = 45194
$ws.Range("E101").Value = 5
$ws.Range("F101").Value = "Fruta"
$ws.Range("G101").Value = 100108
$ws.Range("H101").Value = "Tropicales y subtropicales"
$ws.Range("I101").Value = 100108004
$ws.Range("J101").Value = "Papaya"
$ws.Range("K101").Value = "Cultivar IV Región"
$ws.Range("L101").Value = "Segunda"
$ws.Range("M101").Value = 50
$ws.Range("N101").Value = 15000
$ws.Range("O101").Value = 15000
$ws.Range("P101").Value = 15000
$ws.Range("Q101").Value = "`$/bandeja 10 kilos"
$ws.Range("R101").Value = "Provincia del Elquí"
$ws.Range("S101").Value = 1500
$ws.Range("T101").Value = 10

# Match the date style used by the rest of column D (style index 2 / numFmt 165)
$ws.Range("D100:D101").NumberFormat = $ws.Range("D99").NumberFormat
